$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (C) column for all existing data rows (2-397)
#    from 45181 to 45182.
$ws.Range("C2:C397").Value = 45182

# 2. Row 397 gains an explicit row height (ht="15" customHeight="1"),
#    matching the new rows appended below it.
$ws.Rows.Item(397).RowHeight = 15

# 3. Append the two new records as rows 398 and 399.

# Row 398
$ws.Cells.Item(398, 1).Value = "A 42513-2023"
$ws.Cells.Item(398, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(398, 2).Value = 45180
$ws.Cells.Item(398, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(398, 3).Value = 45182
$ws.Cells.Item(398, 4).Value = "VÄSTERNORRLANDS LÄN"
$ws.Cells.Item(398, 5).Value = "TIMRÅ"
$ws.Cells.Item(398, 6).Value = "SCA"
$ws.Cells.Item(398, 7).Value = 1.4
$ws.Cells.Item(398, 8).Value = 0
$ws.Cells.Item(398, 9).Value = 0
$ws.Cells.Item(398, 10).Value = 0
$ws.Cells.Item(398, 11).Value = 0
$ws.Cells.Item(398, 12).Value = 0
$ws.Cells.Item(398, 13).Value = 0
$ws.Cells.Item(398, 14).Value = 0
$ws.Cells.Item(398, 15).Value = 0
$ws.Cells.Item(398, 16).Value = 0
$ws.Cells.Item(398, 17).Value = 0
$ws.Cells.Item(398, 18).WrapText = $true
$ws.Rows.Item(398).RowHeight = 15

# Row 399
$ws.Cells.Item(399, 1).Value = "A 42516-2023"
$ws.Cells.Item(399, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(399, 2).Value = 45180
$ws.Cells.Item(399, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(399, 3).Value = 45182
$ws.Cells.Item(399, 4).Value = "VÄSTERNORRLANDS LÄN"
$ws.Cells.Item(399, 5).Value = "TIMRÅ"
$ws.Cells.Item(399, 6).Value = "SCA"
$ws.Cells.Item(399, 7).Value = 2.7
$ws.Cells.Item(399, 8).Value = 0
$ws.Cells.Item(399, 9).Value = 0
$ws.Cells.Item(399, 10).Value = 0
$ws.Cells.Item(399, 11).Value = 0
$ws.Cells.Item(399, 12).Value = 0
$ws.Cells.Item(399, 13).Value = 0
$ws.Cells.Item(399, 14).Value = 0
$ws.Cells.Item(399, 15).Value = 0
$ws.Cells.Item(399, 16).Value = 0
$ws.Cells.Item(399, 17).Value = 0
$ws.Cells.Item(399, 18).WrapText = $true
